$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of cell reference -> new text value (values are stored as plain text,
# matching the workbook's existing inlineStr/text cells for this column).
$updates = @{
    "D2" = "330.01"
    "E2" = "7.14%"
    "D3" = "40.11"
    "E3" = "7.86%"
    "D4" = "5.383"
    "D5" = "0.08103"
    "E5" = "3.50%"
    "D6" = "4.527"
    "E6" = "2.93%"
    "D7" = "8.652"
    "E7" = "4.88%"
    "D8" = "1.928"
    "E8" = "1.96%"
    "E9" = "0.41%"
    "D10" = "0.9430"
    "E10" = "2.31%"
    "D11" = "0.1362"
    "E11" = "25.93%"
    "D12" = "0.1976"
    "E12" = "4.14%"
    "D13" = "0.09291"
    "E13" = "4.25%"
    "D14" = "0.03565"
    "E14" = "6.52%"
    "D15" = "0.09592"
    "E15" = "0.22%"
    "D16" = "0.001330"
    "E16" = "-3.40%"
    "D17" = "0.006193"
    "E17" = "9.26%"
    "E18" = "-1.40%"
    "E19" = "2.83%"
    "D20" = "7.231"
    "E20" = "15.13%"
    "E21" = "3.49%"
    "E22" = "4.40%"
    "E23" = "1.70%"
    "D24" = "0.001221"
    "E24" = "2.35%"
    "E25" = "0.23%"
    "E26" = "-14.22%"
    "E39" = "14.35%"
    "D40" = "0.05216"
    "E40" = "3.53%"
    "D41" = "0.007578"
    "E41" = "0.65%"
    "D42" = "0.1428"
    "E42" = "5.77%"
    "D43" = "0.009168"
    "E43" = "5.75%"
    "E44" = "6.41%"
    "D45" = "0.01079"
    "E45" = "36.78%"
    "D46" = "0.00006581"
    "E46" = "0.41%"
    "E47" = "0.07%"
    "E48" = "139.47%"
    "E50" = "0.07%"
    "E51" = "0.07%"
}

foreach ($cellRef in $updates.Keys) {
    $cell = $ws.Range($cellRef)
    # Force text interpretation so values like "330.01" or "7.14%" are not
    # auto-converted to numbers/percent by Excel's smart input parsing.
    $cell.NumberFormat = "@"
    $cell.Value = $updates[$cellRef]
    # Reset the style back to the workbook default (Normal) so the only
    # observable change is the cell text, matching the source diff.
    $cell.Style = "Normal"
}
